$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; existing rows 9-21 shift down to 10-22.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new data record.
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C9").Value = "Los Lagos"
$ws.Range("D9").Value = 44810
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 110
$ws.Range("K9").Value = 22000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 22000
$ws.Range("N9").Value = "$/malla 15 kilos"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 1467
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
